# Update crypto price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.602.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.088.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.89%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("E10").Value = "  -0.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.373"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.615.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.23%  "

$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.683.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.087.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "337.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.98%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.54%  "

$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0929"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.24%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("E32").Value = "  -3.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "153.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  -2.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0688"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.124.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.285.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "

$ws.Range("E45").Value = "  +2.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.89%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.948"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("E49").Value = "  -3.77%  "

$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.693"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
